$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Rename the two existing "build" wizard keys to the new "first" step naming;
# the translated text for these two rows stays the same.
$ws.Cells.Item(189, 2).Value = "lab.wizard.build.first.title"
$ws.Cells.Item(190, 2).Value = "lab.wizard.build.first.subtitle"

# New rows for the wizard tabs (first/coil/atomizer/cotton/build), each pair
# being a "<key>.label" / "<key>.description" translation row.
$newRows = @(
    @{ Row = 191; Key = "lab.wizard.build.first.tab.label"; Value = "Úvod" },
    @{ Row = 192; Key = "lab.wizard.build.first.tab.description"; Value = "Průvodce novým buildem" },
    @{ Row = 193; Key = "lab.wizard.build.coil.tab.label"; Value = "Spirálka" },
    @{ Row = 194; Key = "lab.wizard.build.coil.tab.description"; Value = "Vyberte prosím použitou spirálku" },
    @{ Row = 195; Key = "lab.wizard.build.atomizer.tab.label"; Value = "Atomizér" },
    @{ Row = 196; Key = "lab.wizard.build.atomizer.tab.description"; Value = "Vyberte prosím použitý atomizér" },
    @{ Row = 197; Key = "lab.wizard.build.cotton.tab.label"; Value = "Vata" },
    @{ Row = 198; Key = "lab.wizard.build.cotton.tab.description"; Value = "Vyberte prosím vatu" },
    @{ Row = 199; Key = "lab.wizard.build.build.tab.label"; Value = "Build" },
    @{ Row = 200; Key = "lab.wizard.build.build.tab.description"; Value = "Doplňující informace o buildu" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $rowRange = $ws.Range("A$r`:C$r")
    $rowRange.WrapText = $true
    $rowRange.Font.Size = 10
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $entry.Key
    $ws.Cells.Item($r, 3).Value = $entry.Value
}

$ws.Activate()
$ws.Range("B193").Select()
